$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = '2024-10-05'
$ws.Range("C2").Value = '南宁·花海演绎二次元水上派对'
$ws.Range("D2").Value = '月湾路凤岭儿童公园北门东侧约70米 凤岭儿童公园'
$ws.Range("E2").Value = '2024.10.05 14:00-10.05 21:00'
$ws.Range("F2").Value = 189
$ws.Range("G2").Value = 55
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=92559'
$ws.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202409/MTS1nDly1726642819177.png'
$ws.Range("B3").Value = '2024-10-19'
$ws.Range("C3").Value = '南宁·10.19剑网3同人only——寒光烈火·阵营PK战'
$ws.Range("D3").Value = '大学东路158号 维也纳酒店动物园店'
$ws.Range("E3").Value = '2024.10.19 10:00-10.19 17:30'
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 78
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=92730'
$ws.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202409/3NbN422C1726052875488.jpeg'
$ws.Range("B4").Value = '2024-10-26'
$ws.Range("C4").Value = '南宁·熊喵M动漫嘉年华·万圣派对'
$ws.Range("D4").Value = '亭洪路45号 百益上河城'
$ws.Range("E4").Value = '2024.10.26 11:00-10.27 21:00'
$ws.Range("F4").Value = 109
$ws.Range("G4").Value = 60
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=91894'
$ws.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202409/hBNwkgri1727595368161.jpeg'
$ws.Range("B5").Value = '2024-11-02'
$ws.Range("C5").Value = '南宁·万圣漫控嘉年华10'
$ws.Range("D5").Value = '亭洪路45号 百益上河城'
$ws.Range("E5").Value = '2024.11.02 11:00-11.03 22:00'
$ws.Range("F5").Value = 630
$ws.Range("G5").Value = 50
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202409/mDTW4lHL1727235917704.jpeg'
$ws.Range("B6").Value = '2024-11-02'
$ws.Range("C6").Value = '南宁·梦中礼Lolita茶会'
$ws.Range("D6").Value = '吉兴西路盛天汇一、三、四层 云庭汇·安吉宴会厅'
$ws.Range("E6").Value = '2024.11.02 13:00-11.02 17:00'
$ws.Range("F6").Value = 53
$ws.Range("G6").Value = 138
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=92826'
$ws.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202409/09AXaAJA1726816540668.jpeg'
$ws.Rows.Item(7).Delete()

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = '2024-10-19'
$ws.Range("C2").Value = '南宁·井草圣二 2024《落叶轻扬》指弹吉他音乐会'
$ws.Range("D2").Value = '亭洪路45号 上河城艺术中心'
$ws.Range("E2").Value = '2024.10.19 19:30-10.19 21:00'
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 260
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=91345'
$ws.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202408/7rcuyrqP1724741707788.jpeg'
$ws.Rows.Item(3).Delete()

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = '2024-10-05'
$ws.Range("C2").Value = '南宁·花海演绎二次元水上派对'
$ws.Range("D2").Value = '月湾路凤岭儿童公园北门东侧约70米 凤岭儿童公园'
$ws.Range("E2").Value = '2024.10.05 14:00-10.05 21:00'
$ws.Range("F2").Value = 189
$ws.Range("G2").Value = 55
$ws.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=92559'
$ws.Range("I2").Value = '//i0.hdslb.com/bfs/openplatform/202409/MTS1nDly1726642819177.png'
$ws.Range("B3").Value = '2024-10-19'
$ws.Range("C3").Value = '南宁·10.19剑网3同人only——寒光烈火·阵营PK战'
$ws.Range("D3").Value = '大学东路158号 维也纳酒店动物园店'
$ws.Range("E3").Value = '2024.10.19 10:00-10.19 17:30'
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 78
$ws.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=92730'
$ws.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202409/3NbN422C1726052875488.jpeg'
$ws.Range("B4").Value = '2024-10-19'
$ws.Range("C4").Value = '南宁·井草圣二 2024《落叶轻扬》指弹吉他音乐会'
$ws.Range("D4").Value = '亭洪路45号 上河城艺术中心'
$ws.Range("E4").Value = '2024.10.19 19:30-10.19 21:00'
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 260
$ws.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=91345'
$ws.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202408/7rcuyrqP1724741707788.jpeg'
$ws.Range("B5").Value = '2024-10-26'
$ws.Range("C5").Value = '南宁·熊喵M动漫嘉年华·万圣派对'
$ws.Range("D5").Value = '亭洪路45号 百益上河城'
$ws.Range("E5").Value = '2024.10.26 11:00-10.27 21:00'
$ws.Range("F5").Value = 109
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=91894'
$ws.Range("I5").Value = '//i2.hdslb.com/bfs/openplatform/202409/hBNwkgri1727595368161.jpeg'
$ws.Range("B6").Value = '2024-11-02'
$ws.Range("C6").Value = '南宁·万圣漫控嘉年华10'
$ws.Range("D6").Value = '亭洪路45号 百益上河城'
$ws.Range("E6").Value = '2024.11.02 11:00-11.03 22:00'
$ws.Range("F6").Value = 630
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=87820'
$ws.Range("I6").Value = '//i0.hdslb.com/bfs/openplatform/202409/mDTW4lHL1727235917704.jpeg'
$ws.Range("B7").Value = '2024-11-02'
$ws.Range("C7").Value = '南宁·梦中礼Lolita茶会'
$ws.Range("D7").Value = '吉兴西路盛天汇一、三、四层 云庭汇·安吉宴会厅'
$ws.Range("E7").Value = '2024.11.02 13:00-11.02 17:00'
$ws.Range("F7").Value = 53
$ws.Range("G7").Value = 138
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=92826'
$ws.Range("I7").Value = '//i2.hdslb.com/bfs/openplatform/202409/09AXaAJA1726816540668.jpeg'
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
